$d = $word.ActiveDocument

# Locate the single paragraph whose full text matches the given wildcard
# pattern (the diff only touches a handful of specific paragraphs, each
# uniquely identifiable by a distinctive substring of their text).
function Get-ParaByText($doc, $matchText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $matchText) {
            return $p.Range
        }
    }
    return $null
}

# 1) 'Để chạy thành công mã nguồn mở cần bao nhiêu thời gian' -> split off '2 tuần'
$p1 = Get-ParaByText $d '*Để chạy thành công*cần bao nhiêu thời gian*'
if ($p1 -eq $null) { throw 'paragraph 1 not found' }
$p1.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4CAA1D6E" w14:textId="749A968B" w:rsidR="00612FB1" w:rsidRDefault="000C6EB0" w:rsidP="00017E86"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Để chạy thành công </w:t></w:r><w:r w:rsidR="00017E86"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>mã nguồn</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> mở</w:t></w:r><w:r w:rsidR="00017E86"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> cần </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>2 tuần</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 2) 'Để hiểu rõ mã nguồn mở cần bao nhiêu thời gian' -> split off '5 ngày'
$p2 = Get-ParaByText $d '*Để hiểu rõ mã nguồn mở cần bao nhiêu thời gian*'
if ($p2 -eq $null) { throw 'paragraph 2 not found' }
$p2.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="698E473B" w14:textId="480537F6" w:rsidR="00327CC3" w:rsidRDefault="00327CC3" w:rsidP="00327CC3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Để hiểu rõ mã nguồn mở cần </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>5 ngày</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 3) 'Để thay đổi giao diện, để chỉnh sửa tính năng cần bao nhiêu thời gian' -> split off '10 ngày' + _GoBack bookmark
$p3 = Get-ParaByText $d '*Để thay đổi giao diện, để chỉnh sửa tính năng*' $null
if ($p3 -eq $null) { throw 'paragraph 3 not found' }
$p3.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0B40C28B" w14:textId="21F2549C" w:rsidR="000123C4" w:rsidRPr="00327CC3" w:rsidRDefault="000123C4" w:rsidP="00327CC3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Để thay đổi giao diện, để chỉnh sửa tính năng </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>10 ngày</w:t></w:r><w:bookmarkStart w:id="15" w:name="_GoBack"/><w:bookmarkEnd w:id="15"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 4) remove the old _GoBack bookmark after 'hạn chế các file có nhiều trang tính'
$p4 = Get-ParaByText $d '*hạn chế các file có nhiều trang tính*' $null
if ($p4 -eq $null) { throw 'paragraph 4 not found' }
$p4.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="342B0712" w14:textId="0ADACE78" w:rsidR="00A43D43" w:rsidRPr="007E4E75" w:rsidRDefault="00A43D43" w:rsidP="00A43D43"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r w:rsidRPr="007E4E75"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Giải pháp</w:t></w:r><w:r w:rsidR="007E4E75" w:rsidRPr="007E4E75"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> xử lý</w:t></w:r><w:r w:rsidR="0072012D"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="001377C7"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>hạn chế các file có nhiều trang tính</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output 'done'
